$d = $word.ActiveDocument

$replacements = @(
    @("73×44=3212", "90×90=8100"),
    @("13×98=1274", "37×62=2294"),
    @("15×56=840",  "21×96=2016"),
    @("43×92=3956", "77×67=5159"),
    @("96×96=9216", "23×87=2001"),
    @("20×46=920",  "31×44=1364"),
    @("27×37=999",  "16×68=1088"),
    @("12×38=456",  "80×96=7680"),
    @("85×15=1275", "79×42=3318"),
    @("39×77=3003", "11×43=473"),
    @("32×52=1664", "97×80=7760"),
    @("91×48=4368", "33×40=1320"),
    @("22×97=2134", "35×79=2765"),
    @("78×62=4836", "26×36=936"),
    @("24×82=1968", "60×69=4140"),
    @("78×50=3900", "39×46=1794"),
    @("50×31=1550", "20×11=220"),
    @("57×98=5586", "39×51=1989"),
    @("75×55=4125", "75×32=2400"),
    @("19×27=513",  "80×32=2560"),
    @("90×34=3060", "45×97=4365"),
    @("45×46=2070", "44×72=3168"),
    @("20×43=860",  "44×27=1188"),
    @("33×97=3201", "44×76=3344"),
    @("52×24=1248", "66×75=4950")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
